$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# People : Last Name
$ws.Range("B3").Value = "last_name"
$ws.Range("C3").Value = "People's Last Name"
$ws.Range("D3").Value = "varchar"
$ws.Range("E3").Value = 25
$ws.Range("F3").Value = "NOT NULL"

# People : First Name
$ws.Range("B4").Value = "first_name"
$ws.Range("C4").Value = "People's First Name"
$ws.Range("D4").Value = "varchar"
$ws.Range("E4").Value = 25
$ws.Range("F4").Value = "NOT NULL"

# Customer : birth_date
$ws.Range("A5").Value = "Customer :"
$ws.Range("B5").Value = "birth_date"
$ws.Range("C5").Value = "Customer's birth date"
$ws.Range("D5").Value = "date"
$ws.Range("F5").Value = "NOT NULL"

# Customer : first_order_date
$ws.Range("B6").Value = "first_order_date"
$ws.Range("C6").Value = "Customer's first order date"
$ws.Range("D6").Value = "date"
$ws.Range("F6").Value = "NOT NULL"

# Staff : hiring_date
$ws.Range("A7").Value = "Staff :"
$ws.Range("B7").Value = "hiring_date"
$ws.Range("C7").Value = "Staff's hiring date"
$ws.Range("D7").Value = "date"
$ws.Range("F7").Value = "NOT NULL"

# Address : id
$ws.Range("A8").Value = "Address :"
$ws.Range("B8").Value = "id"
$ws.Range("F8").Value = "Primary Key - Identity (1,1)"

# Address : last_name
$ws.Range("B9").Value = "last_name"
$ws.Range("F9").Value = "NOT NULL"

# Address : first_name
$ws.Range("B10").Value = "first_name"
$ws.Range("F10").Value = "NOT NULL"

# Address : text
$ws.Range("B11").Value = "text"
$ws.Range("D11").Value = "varchar"
$ws.Range("E11").Value = 225
$ws.Range("F11").Value = "NOT NULL"

# Address : postal_code
$ws.Range("B12").Value = "postal_code"
$ws.Range("D12").Value = "varchar"
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = "NOT NULL"

# Address : city
$ws.Range("B13").Value = "city"
$ws.Range("D13").Value = "varchar"
$ws.Range("E13").Value = 225
$ws.Range("F13").Value = "NOT NULL"

# Resize columns B and C to fit the new content, matching the rest of the sheet
$ws.Columns.Item(2).ColumnWidth = 14.17
$ws.Columns.Item(3).ColumnWidth = 23.27

# Reproduce the final cell selection left by the author
$ws.Range("B15").Select()
